$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.274
$ws.Range("D3").Value = 0.274
$ws.Range("G2").Value = -0.01452784503631961
$ws.Range("G3").Value = -0.01452784503631961
$ws.Range("H2").Value = -0.01452784503631961
$ws.Range("H3").Value = -0.01452784503631961
$ws.Range("I2").Value = -0.0188861985472155
$ws.Range("I3").Value = -0.0188861985472155
$ws.Range("J2").Value = -0.0188861985472155
$ws.Range("J3").Value = -0.0188861985472155
$ws.Range("K2").Value = -36.1
$ws.Range("K3").Value = -36.1
$ws.Range("L2").Value = -0.4370460048426151
$ws.Range("L3").Value = -0.4370460048426151
$ws.Range("U2").Value = 9.23
$ws.Range("U3").Value = 9.23
$ws.Range("V2").Value = 0.8315315315315316
$ws.Range("V3").Value = 0.8315315315315316
$ws.Range("W2").Value = -1.058651026392962
$ws.Range("W3").Value = -1.058651026392962
$ws.Range("X2").Value = 0.07516306640179389
$ws.Range("X3").Value = 0.07516306640179389
$ws.Range("Y2").Value = -1.133814092794756
$ws.Range("Y3").Value = -1.133814092794756
$ws.Range("Z2").Value = 5.816901408450704
$ws.Range("Z3").Value = 5.816901408450704
$ws.Range("AA2").Value = -0.1098591549295775
$ws.Range("AA3").Value = -0.1098591549295775
$ws.Range("AB2").Value = 0.0445146645977909
$ws.Range("AB3").Value = 0.0445146645977909
$ws.Range("AC2").Value = -0.1543738195273684
$ws.Range("AC3").Value = -0.1543738195273684
$ws.Range("AD2").Value = 16.1
$ws.Range("AD3").Value = 16.1
$ws.Range("AF2").Value = 16.1
$ws.Range("AF3").Value = 16.1
$ws.Range("AG2").Value = 6.870000000000001
$ws.Range("AG3").Value = 6.870000000000001
$ws.Range("AH2").Value = 0.5919117647058824
$ws.Range("AH3").Value = 0.5919117647058824
$ws.Range("AI2").Value = 0.6686046511627907
$ws.Range("AI3").Value = 0.6686046511627907
$ws.Range("AJ2").Value = 0.3823038397328882
$ws.Range("AJ3").Value = 0.3823038397328882
$ws.Range("AK2").Value = 0.4626262626262627
$ws.Range("AK3").Value = 0.4626262626262627
$ws.Range("AL2").Value = 1.66
$ws.Range("AL3").Value = 1.66
$ws.Range("AM2").Value = 1.66
$ws.Range("AM3").Value = 1.66
$ws.Range("AN2").Value = -38.98305084745763
$ws.Range("AN3").Value = -38.98305084745763
$ws.Range("AO2").Value = -0.9397590361445783
$ws.Range("AO3").Value = -0.9397590361445783
$ws.Range("AP2").Value = -16.63438256658596
$ws.Range("AP3").Value = -16.63438256658596
$ws.Range("AQ2").Value = -0.9397590361445783
$ws.Range("AQ3").Value = -0.9397590361445783
